$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Cell A7 holds a rich-text string: italic "Ca. " + "Chl." + " canadensis L304-6D"
# (species name correction: canadensis -> canadense)
$cell = $ws.Range("A7")
$cell.Value = "Ca. Chl. canadense L304-6D"

# Restore the run-level formatting that setting .Value resets:
# run 1, chars 1-4 "Ca. " -> italic Times New Roman 12pt black
$run1 = $cell.Characters(1, 4)
$run1.Font.Italic = $true
$run1.Font.Size = 12
$run1.Font.Color = 0
$run1.Font.Name = "Times New Roman"

# run 2, chars 5-26 "Chl. canadense L304-6D" -> regular Times New Roman 12pt black
$run2 = $cell.Characters(5, 22)
$run2.Font.Italic = $false
$run2.Font.Size = 12
$run2.Font.Color = 0
$run2.Font.Name = "Times New Roman"

# The column was sized with "best fit"; now that the text is shorter, re-fit it
# to the new (narrower) content width.
$ws.Columns.Item(1).ColumnWidth = 24.75
